$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ParaAfter {
    param($anchorPara, [string]$innerXml)
    $anchorPara.Range.InsertParagraphAfter()
    $newIndex = $anchorPara.Index + 1
    $newPara = $d.Paragraphs.Item($newIndex)
    $fullXml = "<w:p $wns>" + $innerXml + "</w:p>"
    $newPara.Range.InsertXML($fullXml)
    return $d.Paragraphs.Item($newIndex)
}

# --- Hunk 2: insert after 'The temperature of the material...' paragraph (done first to keep indices stable) ---
$anchor2 = $d.Paragraphs.Item(50)
if ($anchor2.Range.Text.Trim() -ne "The temperature of the material should be over 1,000 Kelvin to produce enough electrons.") {
    throw "Anchor2 paragraph text mismatch: $($anchor2.Range.Text)"
}
$h2_0 = Insert-ParaAfter $anchor2 ''
$h2_1 = Insert-ParaAfter $h2_0 '<w:r><w:t xml:space="preserve">From the </w:t></w:r><w:r><w:t>Goebel_06_Chap6_cathodes.pdf</w:t></w:r><w:r><w:t>, it seems we are using the Type A configuration shown on page 6 of the pdf. So it “operates at lower currents and relatively high internal gas pressures, and are heated by orifice heating. “</w:t></w:r>'
$h2_2 = Insert-ParaAfter $h2_1 ''
$h2_3 = Insert-ParaAfter $h2_2 '<w:r><w:t xml:space="preserve">Page 5 of </w:t></w:r><w:r><w:t>Goebel_06_Chap6_cathodes.pdf</w:t></w:r><w:r><w:t xml:space="preserve"> details the three types of self-heating mechanisms: orifice heating, ion heating, election heating. </w:t></w:r>'

# --- Hunk 1: insert after the Goebel_06_Chap6_cathodes.pdf link paragraph (paragraph 15 in original numbering) ---
$anchor1 = $d.Paragraphs.Item(15)
if ($anchor1.Range.Text.Trim() -ne "https://descanso.jpl.nasa.gov/SciTechBook/series1/Goebel_06_Chap6_cathodes.pdf") {
    throw "Anchor1 paragraph text mismatch: $($anchor1.Range.Text)"
}
$h1_0 = Insert-ParaAfter $anchor1 '<w:r><w:t>https://descanso.jpl.nasa.gov/SciTechBook/series1/Goebel__cmprsd_opt.pdf</w:t></w:r>'
$h1_1 = Insert-ParaAfter $h1_0 '<w:r><w:t>First is just chapter 6 the other being the full book</w:t></w:r><w:r><w:t xml:space="preserve"> that is found in the power point by Dr. </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Frieman</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r>'
$h1_2 = Insert-ParaAfter $h1_1 ''
$h1_3 = Insert-ParaAfter $h1_2 '<w:r><w:t xml:space="preserve">Two ways to show: </w:t></w:r>'
$h1_4 = Insert-ParaAfter $h1_3 '<w:r><w:t xml:space="preserve">One is showing where the cathode is in the spacecraft, with a close-up showing info about each part of the cathode, </w:t></w:r><w:r><w:t>https://mars.nasa.gov/mars2020/spacecraft/rover/</w:t></w:r><w:r><w:t xml:space="preserve">, similar to this link for the Mars 2020 rover.  </w:t></w:r>'
$h1_5 = Insert-ParaAfter $h1_4 ''
$h1_6 = Insert-ParaAfter $h1_5 ''

Write-Output "Done. Paragraph count = $($d.Paragraphs.Count)"
